$wb = $excel.ActiveWorkbook

# --- Summary sheet: Total Trades 15 -> 16, Win Rate % 40 -> 37.5 ---
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B6").Value = 16
$wsSummary.Range("B9").Value = 37.5

# --- Strategy Status sheet: MarketMaking row Trades 15 -> 16, Win Rate % 40 -> 37.5 ---
$wsStatus = $wb.Worksheets.Item("Strategy Status")
$wsStatus.Range("D4").Value = 16
$wsStatus.Range("G4").Value = 37.5

# --- Append trade #16 (row 17) to both "All Trades" and "MarketMaking" sheets ---
$sheetNames = @("All Trades", "MarketMaking")
foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("A17").Value = 16

    # Date/time columns must stay plain text (not auto-converted to date/time
    # serials), matching the existing inlineStr cells above them.
    $ws.Range("B17").NumberFormat = "@"
    $ws.Range("B17").Value = "2026-02-17"
    $ws.Range("B17").ClearFormats()

    $ws.Range("C17").Value = "12:28:23"

    $ws.Range("D17").Value = "MarketMaking"
    $ws.Range("E17").Value = "UP"
    $ws.Range("F17").Value = 0.07000000000000001
    $ws.Range("G17").Value = 0.06862699999999999
    $ws.Range("H17").Value = "CLOSED"
    $ws.Range("I17").Value = -1.9608
    $ws.Range("J17").Value = -0
    $ws.Range("K17").Value = 100.01
    $ws.Range("L17").Value = 0
    $ws.Range("M17").Value = 0
    $ws.Range("N17").Value = 0.6
    $ws.Range("O17").Value = "Normal spread capture: 19600 bps"
    $ws.Range("P17").Value = "early_exit"
    $ws.Range("Q17").Value = 0.1
}
